$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Revision")

# Update the barrel jack part number to the Arduino-compatible one.
$ws.Range("B8").Value = "CP-202AH-ND"

# Update the selected cell as recorded in the saved view state.
$ws.Activate()
$ws.Range("D5").Select()
